$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "2026-01-08 23:21:28 EST"
$ws.Range("B10").Value = "c1b327ad-ea0c-473d-89e7-ac6a496a5767"
$ws.Range("C10").Value = -215000
$ws.Range("D10").Value = 1806330
$ws.Range("E10").Value = 1591330
$ws.Range("F10").Value = "No"
